# Applies the cryptos-list price/volume refresh described in the commit.
# Numeric-looking Price (column D) values are written with a leading
# apostrophe so Excel keeps them as text (matching the workbook's
# existing inlineStr / text storage for that column) instead of silently
# re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.607.39"
$ws.Range("E2").Value = "  -2.25%  "

$ws.Range("D3").Value = "1.586.10"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'210.92"
$ws.Range("E5").Value = "  -2.32%  "

$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -2.78%  "

$ws.Range("E9").Value = "  -1.28%  "

$ws.Range("E10").Value = "  -2.93%  "

$ws.Range("D11").Value = "'0.0835"
$ws.Range("E11").Value = "  -1.58%  "

$ws.Range("D12").Value = "1.807.93"
$ws.Range("E12").Value = "  -2.70%  "

$ws.Range("D13").Value = "1.596.41"
$ws.Range("E13").Value = "  -2.28%  "

$ws.Range("D14").Value = "'4.04"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("E15").Value = "  -2.91%  "

$ws.Range("D16").Value = "'64.86"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "26.587.80"
$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("D20").Value = "'207.15"
$ws.Range("E20").Value = "  -4.36%  "

$ws.Range("E21").Value = "  -2.85%  "

$ws.Range("E22").Value = "  -3.38%  "

$ws.Range("E23").Value = "  -4.17%  "

$ws.Range("D24").Value = "'8.88"
$ws.Range("E24").Value = "  -2.61%  "

$ws.Range("D25").Value = "'147.16"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'7.35"
$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("E28").Value = "  -3.28%  "

$ws.Range("D29").Value = "'15.27"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("E31").Value = "  -2.10%  "

$ws.Range("E32").Value = "  -4.37%  "

$ws.Range("D33").Value = "'0.661"
$ws.Range("E33").Value = "  +22.38%  "

$ws.Range("D34").Value = "1.330.46"
$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("E35").Value = "  -2.93%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'1.51"
$ws.Range("E36").Value = "  -3.23%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.41"
$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("E38").Value = "  -1.35%  "

$ws.Range("E39").Value = "  -2.26%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "'5.36"
$ws.Range("E41").Value = "  +3.68%  "

$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("E43").Value = "  -3.59%  "

$ws.Range("D44").Value = "'63.40"
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "1.721.51"
$ws.Range("E45").Value = "  -2.53%  "

$ws.Range("D46").Value = "'89.85"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("E48").Value = "  +2.43%  "

$ws.Range("E49").Value = "  -1.82%  "

$ws.Range("D50").Value = "'0.0982"
$ws.Range("E50").Value = "  +2.74%  "

$ws.Range("D51").Value = "'7.48"
$ws.Range("E51").Value = "  -0.65%  "
